$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 181, shifting all existing rows (previously 181-214)
# down to 182-215.
$ws.Rows.Item(181).Insert()

# Fill in the data for the newly inserted row 181 (same shape as the other
# "Feria Lagunitas de Puerto Montt" / Perejil rows in this sheet).
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44617
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112044
$ws.Range("G181").Value = "Perejil"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 180
$ws.Range("K181").Value = 5000
$ws.Range("L181").Value = 5000
$ws.Range("M181").Value = 5000
$ws.Range("N181").Value = "$/docena de atados (3 kilos)"
$ws.Range("O181").Value = "Región Metropolitana"
$ws.Range("P181").Value = 1667
$ws.Range("Q181").Value = 3
$ws.Range("R181").Value = "Hortaliza"
